$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for Brócoli at Terminal
# Hortofrutícola Agro Chillán. It belongs chronologically where row 167
# used to sit, so insert a fresh row there (pushing every following row
# down by one, which is exactly what the diff shows: row 168 now holds
# what used to be row 167's data, row 169 what used to be row 168's, …,
# and a brand-new row 281 holds what used to be the old last row, 280).
$ws.Rows.Item(167).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(167, 1).Value = 7
$ws.Cells.Item(167, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(167, 3).Value = "Ñuble"
$ws.Cells.Item(167, 4).Value = 44762
$ws.Cells.Item(167, 5).Value = 16
$ws.Cells.Item(167, 6).Value = 100112023
$ws.Cells.Item(167, 7).Value = "Brócoli"
$ws.Cells.Item(167, 8).Value = "Sin especificar"
$ws.Cells.Item(167, 9).Value = "Primera"
$ws.Cells.Item(167, 10).Value = 200
$ws.Cells.Item(167, 11).Value = 900
$ws.Cells.Item(167, 12).Value = 1000
$ws.Cells.Item(167, 13).Value = 950
$ws.Cells.Item(167, 14).Value = "$/unidad"
$ws.Cells.Item(167, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(167, 16).Value = 950
$ws.Cells.Item(167, 17).Value = 1
$ws.Cells.Item(167, 18).Value = "Hortaliza"
